$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wallet-transaction rows (rows 5-20) appended below the existing
# TransactionID / RetailerID / Amount / Type / Date / Description rows.

$ws.Range("A5").Value = "aa45270b-e681-4c3a-b135-158beb877dd8"
$ws.Range("B5").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "Credit"
$ws.Range("E5").Value = "2025-08-07 13:31:49"
$ws.Range("F5").Value = "Wallet top-up"

$ws.Range("A6").Value = "925a579f-0d74-49b8-9e60-f47984d8ce12"
$ws.Range("B6").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C6").Value = 800
$ws.Range("D6").Value = "Credit"
$ws.Range("E6").Value = "2025-08-07 14:06:13"
$ws.Range("F6").Value = "Wallet top-up"

$ws.Range("A7").Value = "5094161e-b960-4117-b287-9a6496f3c37a"
$ws.Range("B7").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = "Credit"
$ws.Range("E7").Value = "2025-08-07 14:06:52"
$ws.Range("F7").Value = "Wallet top-up"

$ws.Range("A8").Value = "9d9b44b3-78e7-424c-a863-8382583efc88"
$ws.Range("B8").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "Credit"
$ws.Range("E8").Value = "2025-08-07 14:07:03"
$ws.Range("F8").Value = "Wallet top-up"

$ws.Range("A9").Value = "64d4f5fa-1b9a-4637-95db-3897eac5e11f"
$ws.Range("B9").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C9").Value = 881.53
$ws.Range("D9").Value = "Debit"
$ws.Range("E9").Value = "2025-08-07 21:13:22"
$ws.Range("F9").Value = "Payment for order add6407f-8225-46ae-be70-e5a3c9a9b5c7"

$ws.Range("A10").Value = "55564ace-2267-49e3-a7f7-ce79b06efc35"
$ws.Range("B10").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C10").Value = 881.53
$ws.Range("D10").Value = "Debit"
$ws.Range("E10").Value = "2025-08-07 21:14:03"
$ws.Range("F10").Value = "Payment for order add6407f-8225-46ae-be70-e5a3c9a9b5c7"

$ws.Range("A11").Value = "f0e8a2d7-d35d-47ce-a1aa-3dd6a16e5ac2"
$ws.Range("B11").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C11").Value = 300
$ws.Range("D11").Value = "Credit"
$ws.Range("E11").Value = "2025-08-07 22:54:34"
$ws.Range("F11").Value = "Wallet top-up"

$ws.Range("A12").Value = "b17fd9d3-ab4f-4d88-a217-28850186e939"
$ws.Range("B12").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C12").Value = 200
$ws.Range("D12").Value = "Credit"
$ws.Range("E12").Value = "2025-08-07 22:54:55"
$ws.Range("F12").Value = "Wallet top-up"

$ws.Range("A13").Value = "ac6bb0db-7cca-4998-8690-a8ec41d5e67a"
$ws.Range("B13").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C13").Value = 881.53
$ws.Range("D13").Value = "Debit"
$ws.Range("E13").Value = "2025-08-07 22:54:57"
$ws.Range("F13").Value = "Payment for order add6407f-8225-46ae-be70-e5a3c9a9b5c7"

$ws.Range("A14").Value = "20193007-f5e7-4778-a459-0d2c5535e4a5"
$ws.Range("B14").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C14").Value = 1000
$ws.Range("D14").Value = "Credit"
$ws.Range("E14").Value = "2025-08-07 23:00:24"
$ws.Range("F14").Value = "Wallet top-up"

$ws.Range("A15").Value = "5ed4213b-db4d-4152-ab7a-936c78cd97ef"
$ws.Range("B15").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C15").Value = 881.53
$ws.Range("D15").Value = "Debit"
$ws.Range("E15").Value = "2025-08-07 23:00:27"
$ws.Range("F15").Value = "Payment for order add6407f-8225-46ae-be70-e5a3c9a9b5c7"

$ws.Range("A16").Value = "ca15432e-d936-47da-ad2c-e73e395272dd"
$ws.Range("B16").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C16").Value = 73.5
$ws.Range("D16").Value = "Debit"
$ws.Range("E16").Value = "2025-08-07 23:02:17"
$ws.Range("F16").Value = "Payment for order 89bdc2f6-0e22-47a8-b4f2-b7b5696fc495"

$ws.Range("A17").Value = "659e23a0-0c33-4921-bf3d-290f6b875a33"
$ws.Range("B17").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C17").Value = 2000
$ws.Range("D17").Value = "Credit"
$ws.Range("E17").Value = "2025-08-07 23:11:56"
$ws.Range("F17").Value = "Wallet top-up"

$ws.Range("A18").Value = "e3965d49-10ce-4f9a-913c-59988cee7b41"
$ws.Range("B18").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C18").Value = 73.5
$ws.Range("D18").Value = "Debit"
$ws.Range("E18").Value = "2025-08-07 23:11:58"
$ws.Range("F18").Value = "Payment for order 89bdc2f6-0e22-47a8-b4f2-b7b5696fc495"

$ws.Range("A19").Value = "a24054ed-213d-4d55-9f7d-7937065399d5"
$ws.Range("B19").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C19").Value = 881.53
$ws.Range("D19").Value = "Debit"
$ws.Range("E19").Value = "2025-08-07 23:11:58"
$ws.Range("F19").Value = "Payment for order add6407f-8225-46ae-be70-e5a3c9a9b5c7"

$ws.Range("A20").Value = "72f84b03-c371-4829-bb7d-a19d4b4b1144"
$ws.Range("B20").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C20").Value = 182.9
$ws.Range("D20").Value = "Debit"
$ws.Range("E20").Value = "2025-08-07 23:14:24"
$ws.Range("F20").Value = "Payment for order 0947da20-6ab3-444d-97b4-2aa9c1662a75"
